$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for 2025-10-25: write the date as literal text (leading apostrophe
# forces text interpretation, matching how the other date cells in this sheet
# are stored as plain strings rather than date serials), then reset the
# style Excel stamps on quote-prefixed cells so it matches the rest of the
# column (no explicit style).
$ws.Cells.Item(69, 1).Value = "'10/25/2025"
$ws.Cells.Item(69, 1).Style = "Normal"

$ws.Cells.Item(69, 2).Value = 10924.94
